$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.721.85"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.592.26"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'209.01"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'22.31"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "1.594.31"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'3.83"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "27.718.89"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "'218.07"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'4.16"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'9.75"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "'153.84"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'7.06"
$ws.Range("E26").Value = "  +5.49%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'15.13"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "1.382.16"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "'1.55"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "'0.969"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Value = "'0.0169"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Value = "'0.534"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "'0.827"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'0.986"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "'64.49"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("D45").Value = "'5.26"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "1.730.60"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "'86.06"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "'0.0968"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -0.22%  "
